# Generate Report for Handoff
# Adds a new "Ready for handoff" entry (f53f2a61-6a0c-42b7-b117-f7937ecc2f78.md)
# as row 9 to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$fileBase = "f53f2a61-6a0c-42b7-b117-f7937ecc2f78"
$mdName   = "$fileBase.md"
$hoHash   = "68446a203c5de11450dec4fa3e2009c10aebfddb"
$commit   = "68446a203c5de11450dec4fa3e2009c10aebfddb"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = $mdName
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-19 08:42:34"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$mdName",
    "",
    "",
    "e2e\$mdName"
) | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("B9").Value = ".md"
$wsZh.Range("C9").Value = "Ready for handoff"
$wsZh.Range("D9").Value = "e2e"
$wsZh.Range("E9").Value = "ht"
$wsZh.Range("F9").Value = "False"
$wsZh.Range("G9").Value = "$fileBase.$hoHash.zh-cn.xlf"
$wsZh.Range("H9").Value = "2016-08-19 08:42:30"
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I9").Value = ""
$wsZh.Range("J9").Value = ""
$wsZh.Range("K9").Value = "0001-01-01 00:00:00"
$wsZh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L9").Value = ""
$wsZh.Range("M9").Value = "True"
$wsZh.Range("N9").Value = ""
$wsZh.Range("O9").Value = "False"
$wsZh.Range("P9").Value = ""

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("B9").Value = ".md"
$wsDe.Range("C9").Value = "Ready for handoff"
$wsDe.Range("D9").Value = "e2e"
$wsDe.Range("E9").Value = "ht"
$wsDe.Range("F9").Value = "False"
$wsDe.Range("G9").Value = "$fileBase.$hoHash.de-de.xlf"
$wsDe.Range("H9").Value = "2016-08-19 08:42:34"
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I9").Value = ""
$wsDe.Range("J9").Value = ""
$wsDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L9").Value = ""
$wsDe.Range("M9").Value = "True"
$wsDe.Range("N9").Value = ""
$wsDe.Range("O9").Value = "False"
$wsDe.Range("P9").Value = ""

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null
